{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the relevant list-item paragraphs by their current text so the\n// script is resilient to any incidental paragraph-index differences.\nlet computerLogicPara = null;\nlet compShouldPlayPara = null;\nlet oldPhpPara = null;\nlet countUpScorePara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text.trim();\n  if (text === \"Computer logic\") {\n    computerLogicPara = paragraphs.items[i];\n  } else if (text === \"Comp should play better when question is to play anywhere\") {\n    compShouldPlayPara = paragraphs.items[i];\n  } else if (text === \"PHP\") {\n    oldPhpPara = paragraphs.items[i];\n  } else if (text === \"Count up score at end\") {\n    countUpScorePara = paragraphs.items[i];\n  }\n}\n\n// \"Computer logic\" becomes \"PHP\".\nif (computerLogicPara) {\n  computerLogicPara.insertText(\"PHP\", Word.InsertLocation.replace);\n}\n\n// Remove the now-obsolete bullets.\nif (compShouldPlayPara) {\n  compShouldPlayPara.delete();\n}\nif (oldPhpPara) {\n  oldPhpPara.delete();\n}\nif (countUpScorePara) {\n  countUpScorePara.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# \"Computer logic\" becomes \"PHP\".\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Computer logic\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"PHP\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# Collect the (1-based) indices of every paragraph that must be removed:\n#   - \"Comp should play better when question is to play anywhere\"\n#   - \"Count up score at end\"\n#   - the now-duplicate old \"PHP\" bullet (keep only the first \"PHP\" paragraph,\n#     i.e. the one that used to be \"Computer logic\")\n$removeText = @(\n    \"Comp should play better when question is to play anywhere\",\n    \"Count up score at end\"\n)\n\n$indicesToDelete = @()\n$phpIndices = @()\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.Trim()\n    if ($removeText -contains $t) {\n        $indicesToDelete += $i\n    }\n    if ($t -eq \"PHP\") {\n        $phpIndices += $i\n    }\n    $i++\n}\n\nif ($phpIndices.Count -gt 1) {\n    for ($k = 1; $k -lt $phpIndices.Count; $k++) {\n        $indicesToDelete += $phpIndices[$k]\n    }\n}\n\n# Delete from the highest index down so earlier indices stay valid.\n$indicesToDelete = $indicesToDelete | Sort-Object -Descending -Unique\nforeach ($idx in $indicesToDelete) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
